$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# --- helper fragments -------------------------------------------------
# Paragraph 1: "String nombreGato;" -> same text, split into proofed runs
$xmlP1 = "<w:p xmlns:w='$wNs'>" +
         "<w:proofErr w:type='spellStart'/>" +
         "<w:r><w:t>String</w:t></w:r>" +
         "<w:proofErr w:type='spellEnd'/>" +
         "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
         "<w:proofErr w:type='spellStart'/>" +
         "<w:r><w:t>nombreGato</w:t></w:r>" +
         "<w:proofErr w:type='spellEnd'/>" +
         "<w:r><w:t>;</w:t></w:r>" +
         "</w:p>"

# Paragraph 2: "String razaGato;" -> same text, split into proofed runs
$xmlP2 = "<w:p xmlns:w='$wNs'>" +
         "<w:proofErr w:type='spellStart'/>" +
         "<w:r><w:t>String</w:t></w:r>" +
         "<w:proofErr w:type='spellEnd'/>" +
         "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
         "<w:proofErr w:type='spellStart'/>" +
         "<w:r><w:t>razaGato</w:t></w:r>" +
         "<w:proofErr w:type='spellEnd'/>" +
         "<w:r><w:t>;</w:t></w:r>" +
         "</w:p>"

# New paragraph 3 (inserted): "Int uñasGato;" carries the _GoBack bookmark
$xmlP3New = "<w:p xmlns:w='$wNs'>" +
         "<w:proofErr w:type='spellStart'/>" +
         "<w:r><w:t>Int</w:t></w:r>" +
         "<w:proofErr w:type='spellEnd'/>" +
         "<w:r><w:t xml:space='preserve'> u" + [char]0x00F1 + "asGato;</w:t></w:r>" +
         "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" +
         "<w:bookmarkEnd w:id='0'/>" +
         "</w:p>"

# Paragraph 4 (was paragraph 3): "String colorGato;" -> split into proofed runs
$xmlP4 = "<w:p xmlns:w='$wNs'>" +
         "<w:proofErr w:type='spellStart'/>" +
         "<w:r><w:t>String</w:t></w:r>" +
         "<w:proofErr w:type='spellEnd'/>" +
         "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
         "<w:proofErr w:type='spellStart'/>" +
         "<w:r><w:t>colorGato</w:t></w:r>" +
         "<w:proofErr w:type='spellEnd'/>" +
         "<w:r><w:t>;</w:t></w:r>" +
         "</w:p>"

# Paragraph 5 (was paragraph 4): "Int tamañoGato;" -> split into proofed runs,
# bookmark no longer lives here
$xmlP5 = "<w:p xmlns:w='$wNs'>" +
         "<w:proofErr w:type='spellStart'/>" +
         "<w:r><w:t>Int</w:t></w:r>" +
         "<w:proofErr w:type='spellEnd'/>" +
         "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
         "<w:proofErr w:type='spellStart'/>" +
         "<w:r><w:t>tama" + [char]0x00F1 + "oGato</w:t></w:r>" +
         "<w:proofErr w:type='spellEnd'/>" +
         "<w:r><w:t>;</w:t></w:r>" +
         "</w:p>"

# --- apply --------------------------------------------------------------
# Insert a brand new paragraph right before the (current) 3rd paragraph
# ("String colorGato;") to host the new "Int uñasGato;" attribute line.
$thirdPara = $d.Paragraphs.Item(3)
$thirdPara.Range.InsertParagraphBefore()

# Re-fetch paragraphs now that the document has 5 of them and rewrite each
# one's contents (this also moves the _GoBack bookmark onto the new line).
$d.Paragraphs.Item(1).Range.InsertXML($xmlP1)
$d.Paragraphs.Item(2).Range.InsertXML($xmlP2)
$d.Paragraphs.Item(3).Range.InsertXML($xmlP3New)
$d.Paragraphs.Item(4).Range.InsertXML($xmlP4)
$d.Paragraphs.Item(5).Range.InsertXML($xmlP5)
